$wb = $excel.ActiveWorkbook

# New row data (row 53) for each of the 4 worksheets. Each sheet gets a
# new daily log entry appended after the existing last row (52).
$dateVal = 45839.43578703704
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---- Sheet 1: DE_LFT_#1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(53, 1).Value = $dateVal
$ws.Cells.Item(53, 1).NumberFormat = $dateFmt
$ws.Cells.Item(53, 2).Value = "0x01,0x7c"
$ws.Cells.Item(53, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(53, 4).Value = "0x01,0x60"
$ws.Cells.Item(53, 5).Value = "0x14"
$ws.Cells.Item(53, 6).Value = 380
$ws.Cells.Item(53, 7).Value = 7598631275147109 * [Math]::Pow(10, 8)
$ws.Cells.Item(53, 8).Value = 352
$ws.Cells.Item(53, 9).Value = 14

# ---- Sheet 2: DE_LFT_#2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(53, 1).Value = $dateVal
$ws.Cells.Item(53, 1).NumberFormat = $dateFmt
$ws.Cells.Item(53, 2).Value = "0x01,0x7c"
$ws.Cells.Item(53, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(53, 4).Value = "0x01,0x64"
$ws.Cells.Item(53, 5).Value = "0xe"
$ws.Cells.Item(53, 6).Value = 380
$ws.Cells.Item(53, 7).Value = 568432987514711 * [Math]::Pow(10, 9)
$ws.Cells.Item(53, 8).Value = 356
$ws.Cells.Item(53, 9).Value = 14

# ---- Sheet 3: DE_PLT_#1 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(53, 1).Value = $dateVal
$ws.Cells.Item(53, 1).NumberFormat = $dateFmt
$ws.Cells.Item(53, 2).Value = "0x00,0x82"
$ws.Cells.Item(53, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(53, 4).Value = "0x00,0x7E"
$ws.Cells.Item(53, 5).Value = "0x7"
$ws.Cells.Item(53, 6).Value = 130
$ws.Cells.Item(53, 7).Value = 568631262647114 * [Math]::Pow(10, 9)
$ws.Cells.Item(53, 8).Value = 126
$ws.Cells.Item(53, 9).Value = 7

# ---- Sheet 4: DE_PLT_#2 ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(53, 1).Value = $dateVal
$ws.Cells.Item(53, 1).NumberFormat = $dateFmt
$ws.Cells.Item(53, 2).Value = "0x00,0x82"
$ws.Cells.Item(53, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(53, 4).Value = "0x00,0x7E"
$ws.Cells.Item(53, 5).Value = "0x3"
$ws.Cells.Item(53, 6).Value = 130
$ws.Cells.Item(53, 7).Value = 985046333984776 * [Math]::Pow(10, 9)
$ws.Cells.Item(53, 8).Value = 126
$ws.Cells.Item(53, 9).Value = 3
